# Applies the price/volume refresh captured in the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values.
# Values that look like plain numbers are prefixed with a literal leading
# apostrophe (Excel's text-entry marker) so they are stored as text and keep
# their exact original formatting (trailing zeros, etc.) instead of being
# coerced into floating point numbers.
$ws.Range("D2").Value = "62.747.25"
$ws.Range("D3").Value = "3.473.29"
$ws.Range("D5").Value = "'415.22"
$ws.Range("D6").Value = "'130.03"
$ws.Range("D7").Value = "'0.628"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D10").Value = "'0.153"
$ws.Range("D11").Value = "'42.54"
$ws.Range("D12").Value = "'9.83"
$ws.Range("D13").Value = "'0.0000228"
$ws.Range("D14").Value = "4.024.21"
$ws.Range("D16").Value = "'20.59"
$ws.Range("D17").Value = "3.465.15"
$ws.Range("D18").Value = "'12.61"
$ws.Range("D20").Value = "62.705.06"
$ws.Range("D21").Value = "'463.63"
$ws.Range("D22").Value = "'90.47"
$ws.Range("D24").Value = "'13.27"
$ws.Range("D25").Value = "'10.82"
$ws.Range("D26").Value = "'3.31"
$ws.Range("D27").Value = "'33.32"
$ws.Range("D28").Value = "'4.80"
$ws.Range("D29").Value = "'7.60"
$ws.Range("D33").Value = "'0.113"
$ws.Range("D34").Value = "'41.00"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D36").Value = "'58.18"
$ws.Range("D40").Value = "'149.28"
$ws.Range("D41").Value = "'0.324"
$ws.Range("D45").Value = "'4.43"
$ws.Range("D47").Value = "0.0₃0575"
$ws.Range("D48").Value = "'2.38"
$ws.Range("D49").Value = "'16.36"
$ws.Range("D50").Value = "'22.37"
$ws.Range("D51").Value = "'0.141"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  +8.11%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  +5.36%  "
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("E25").Value = "  +15.91%  "
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  +34.97%  "
$ws.Range("E48").Value = "  +10.77%  "
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  -4.85%  "
